$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 599, shifting existing rows 599..645 down to 600..646
$ws.Rows.Item(599).Insert()

# Populate the newly inserted row 599 with the new price record
$ws.Cells.Item(599, 1).Value = 9
$ws.Cells.Item(599, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(599, 3).Value = "Metropolitana"
$ws.Cells.Item(599, 4).Value = 44585
$ws.Cells.Item(599, 5).Value = 13
$ws.Cells.Item(599, 6).Value = "Fruta"
$ws.Cells.Item(599, 7).Value = 100104
$ws.Cells.Item(599, 8).Value = "Frutos de pepita"
$ws.Cells.Item(599, 9).Value = 100104005
$ws.Cells.Item(599, 10).Value = "Pera"
$ws.Cells.Item(599, 11).Value = "Salvador Izquierdo"
$ws.Cells.Item(599, 12).Value = "Primera"
$ws.Cells.Item(599, 13).Value = 25
$ws.Cells.Item(599, 14).Value = 12000
$ws.Cells.Item(599, 15).Value = 12000
$ws.Cells.Item(599, 16).Value = 12000
$ws.Cells.Item(599, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(599, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(599, 19).Value = 667
$ws.Cells.Item(599, 20).Value = 18
